$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.879.64'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '2.974.74'
$ws.Range("E3").Value = '  +1.46%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.42'
$ws.Range("E5").Value = '  -1.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.57'
$ws.Range("E6").Value = '  -3.63%  '
$ws.Range("E7").Value = '  -2.91%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.601'
$ws.Range("E9").Value = '  -5.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.62'
$ws.Range("E10").Value = '  -4.67%  '
$ws.Range("E11").Value = '  +2.48%  '
$ws.Range("E12").Value = '  -3.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.90'
$ws.Range("E13").Value = '  -4.38%  '
$ws.Range("D14").Value = '3.444.23'
$ws.Range("E14").Value = '  +1.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.49'
$ws.Range("E15").Value = '  -5.97%  '
$ws.Range("D16").Value = '2.962.06'
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.986'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").Value = '51.827.61'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.35'
$ws.Range("E20").Value = '  -3.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.36'
$ws.Range("E21").Value = '  -5.31%  '
$ws.Range("D22").Value = '0.0₃0961'
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.94'
$ws.Range("E23").Value = '  -2.99%  '
$ws.Range("E24").Value = '  -3.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.69'
$ws.Range("E25").Value = '  -5.02%  '
$ws.Range("E26").Value = '  -4.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.69'
$ws.Range("E27").Value = '  -1.92%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("E29").Value = '  +3.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.27'
$ws.Range("E30").Value = '  -3.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.23'
$ws.Range("E31").Value = '  +2.40%  '
$ws.Range("E32").Value = '  -5.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.70'
$ws.Range("E33").Value = '  -7.21%  '
$ws.Range("E34").Value = '  +12.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.82'
$ws.Range("E35").Value = '  -2.95%  '
$ws.Range("E36").Value = '  -4.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("E39").Value = '  +2.55%  '
$ws.Range("E40").Value = '  -4.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.32'
$ws.Range("E41").Value = '  -6.45%  '
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '123.24'
$ws.Range("E44").Value = '  +3.30%  '
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("D46").Value = '2.103.23'
$ws.Range("E46").Value = '  -1.87%  '
$ws.Range("E47").Value = '  -5.93%  '
$ws.Range("E48").Value = '  -7.92%  '
$ws.Range("D49").Value = '3.270.80'
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("E50").Value = '  -3.90%  '
$ws.Range("E51").Value = '  -2.53%  '
